# Rewrite each changed paragraph precisely via Range.InsertXML using a single-part
# WordProcessingML package. This lets us set exact run/xml:space structure and
# paragraph properties (or lack thereof) to match the target, rather than relying
# on property setters that cannot fully clear <w:pPr>.
$d = $word.ActiveDocument

$xml3 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">                </w:t></w:r><w:r><w:t>Blackeye or Field Peas</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para3 = $d.Paragraphs.Item(3)
$rng3 = $para3.Range
$rng3.InsertXML($xml3)

$xml5 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Instructions\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para5 = $d.Paragraphs.Item(5)
$rng5 = $para5.Range
$rng5.InsertXML($xml5)

$xml6 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Wash peas. Then blanch snow peas in boiling water for 1 minute, blanch </w:t></w:r><w:r><w:t xml:space="preserve">sugar-snap peas </w:t></w:r><w:r><w:t xml:space="preserve">for </w:t></w:r><w:r><w:t xml:space="preserve">1 1/2 minutes in boiling </w:t></w:r><w:r><w:t xml:space="preserve">water and </w:t></w:r><w:r><w:t>blanch shelled pea</w:t></w:r><w:r><w:t xml:space="preserve">s </w:t></w:r><w:r><w:t xml:space="preserve">for 2 minutes in boiling water. </w:t></w:r><w:r><w:t>\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para6 = $d.Paragraphs.Item(6)
$rng6 = $para6.Range
$rng6.InsertXML($xml6)

$xml8 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Once the time is up, immediately </w:t></w:r><w:r><w:t>immersing</w:t></w:r><w:r><w:t xml:space="preserve"> the peas in an ice bath until cool.</w:t></w:r><w:r><w:t>\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para8 = $d.Paragraphs.Item(8)
$rng8 = $para8.Range
$rng8.InsertXML($xml8)

$xml10 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Place the peas on screens or trays in the dehydrator or oven.</w:t></w:r><w:r><w:t>\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para10 = $d.Paragraphs.Item(10)
$rng10 = $para10.Range
$rng10.InsertXML($xml10)

$xml12 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Set the temperature to 125 degrees Fahrenheit, or the lowest setting your oven allows, and dehydrate for 5 to 13 hours until the peas crisp and become brittle</w:t></w:r><w:r><w:t>.\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para12 = $d.Paragraphs.Item(12)
$rng12 = $para12.Range
$rng12.InsertXML($xml12)

$xml14 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>In</w:t></w:r><w:r><w:t xml:space="preserve"> an oven, it is necessary to stir the peas several times during drying to prevent uneven drying.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para14 = $d.Paragraphs.Item(14)
$rng14 = $para14.Range
$rng14.InsertXML($xml14)

$xml15 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para15 = $d.Paragraphs.Item(15)
$rng15 = $para15.Range
$rng15.InsertXML($xml15)

$xml16 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Shelling peas become wrinkly and hardened when completely dried. When firm pressure is applied to the peas, they should crumble. </w:t></w:r><w:r><w:t>\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para16 = $d.Paragraphs.Item(16)
$rng16 = $para16.Range
$rng16.InsertXML($xml16)

$xml17 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para17 = $d.Paragraphs.Item(17)
$rng17 = $para17.Range
$rng17.InsertXML($xml17)

$xml18 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Store in airtight glass canning jars.\n</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$para18 = $d.Paragraphs.Item(18)
$rng18 = $para18.Range
$rng18 = $d.Range($rng18.Start, $rng18.End - 1)
$rng18.InsertXML($xml18)
